# RDCC-3540 Upload Template Improvements
# - Rename sheet "Case Worker Data" -> "Staff Data"
# - Rename "Area of Work1..8" header columns (L1:S1) to "Service1..8"
# - Add two new header columns: V1 "Task Supervisor", W1 "Case Allocator"
# - Update selection to the newly added last header cell (W1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the first worksheet
$ws.Name = "Staff Data"

# Add the two brand-new trailing header cells first so the new shared
# strings are appended to the shared-string table in source order
# (Task Supervisor, Case Allocator) ahead of the Service* renames.
$ws.Range("V1").Value = "Task Supervisor"
$ws.Range("W1").Value = "Case Allocator"

# Rename "Area of Work1".."Area of Work8" columns to "Service1".."Service8"
$ws.Range("L1").Value = "Service1"
$ws.Range("M1").Value = "Service2"
$ws.Range("N1").Value = "Service3"
$ws.Range("O1").Value = "Service4"
$ws.Range("P1").Value = "Service5"
$ws.Range("Q1").Value = "Service6"
$ws.Range("R1").Value = "Service7"
$ws.Range("S1").Value = "Service8"

# Move/collapse the active selection onto the last newly-added header cell
$ws.Range("W1").Select()
